# Generate Report for Handoff
# This script updates the localization-status workbook to reflect that the
# "b.md" file has now been handed off for zh-cn / de-de, replacing the
# previous "in sync" status with "Ready for handoff" and recording the
# new handoff package (b.63290e5768f688058c7b37413b0a5c26c308f864.*) plus
# an error message explaining the handback file is stale.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/920d4fc47c76312342ad245ff75dc23ce33abfa5/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c4f8bfabd9bfcc3ffc8eb0ca1588532d16773099/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the b.md entry (zh-cn / de-de status columns)
# ---------------------------------------------------------------------
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-02 18:42:01"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the b.md entry
#   C Status, F Content Duplicate, G Latest Handoff File,
#   H Latest Handoff Datetime, P Error Detail
# ---------------------------------------------------------------------
$zhcn.Range("C3").Value = "Ready for handoff"

# "False" must stay a text value (matching the existing "True"/"False"
# strings used elsewhere in the sheet) rather than becoming a native
# Boolean cell, so force text with a leading apostrophe and then clear
# the resulting quote-prefix style.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"

$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-02 18:41:56"
$zhcn.Range("P3").Value = $errorDetail

# Widen the Error Detail column so the long message is readable.
$zhcn.Columns.Item(16).ColumnWidth = 39.14

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the b.md entry
# ---------------------------------------------------------------------
$dede.Range("C3").Value = "Ready for handoff"

$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"

$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-02 18:42:01"
$dede.Range("P3").Value = $errorDetail

$dede.Columns.Item(16).ColumnWidth = 39.14
